$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'27.400.39"
$ws.Range("E2").Value = '  +6.37%  '

$ws.Range("D3").Value = "'1.816.77"
$ws.Range("E3").Value = '  +6.46%  '

$ws.Range("E4").Value = '  +0.13%  '

$ws.Range("D5").Value = "'344.54"
$ws.Range("E5").Value = '  +3.98%  '

$ws.Range("D6").Value = "'0.9999"
$ws.Range("E6").Value = '  +0.20%  '

$ws.Range("D7").Value = "'0.3851"
$ws.Range("E7").Value = '  +4.65%  '

$ws.Range("D8").Value = "'50.43"
$ws.Range("E8").Value = '  +3.88%  '

$ws.Range("D9").Value = "'0.3535"
$ws.Range("E9").Value = '  +7.02%  '

$ws.Range("D10").Value = "'1.242"
$ws.Range("E10").Value = '  +6.41%  '

$ws.Range("D11").Value = "'0.07777"
$ws.Range("E11").Value = '  +6.05%  '

$ws.Range("B12").Value = 'BinanceUSD'
$ws.Range("C12").Value = 'https://coinranking.com/coin/vSo2fu9iE1s0Y+binanceusd-busd'
$ws.Range("D12").Value = "'1.002"
$ws.Range("E12").Value = '  +0.26%  '

$ws.Range("B13").Value = 'Solana'
$ws.Range("C13").Value = 'https://coinranking.com/coin/zNZHO_Sjf+solana-sol'
$ws.Range("D13").Value = "'22.67"
$ws.Range("E13").Value = '  +13.80%  '

$ws.Range("D14").Value = "'6.648"
$ws.Range("E14").Value = '  +7.32%  '

$ws.Range("D15").Value = "'7.255"
$ws.Range("E15").Value = '  +5.77%  '

$ws.Range("D16").Value = "'1.816.47"
$ws.Range("E16").Value = '  +6.99%  '

$ws.Range("D17").Value = "'0.00001127"
$ws.Range("E17").Value = '  +5.32%  '

$ws.Range("D18").Value = "'0.06798"
$ws.Range("E18").Value = '  +2.81%  '

$ws.Range("D19").Value = "'87.34"
$ws.Range("E19").Value = '  +7.59%  '

$ws.Range("D20").Value = "'1.000"
$ws.Range("E20").Value = '  +0.26%  '

$ws.Range("D21").Value = "'17.91"
$ws.Range("E21").Value = '  +10.76%  '

$ws.Range("D22").Value = "'6.569"
$ws.Range("E22").Value = '  +8.57%  '

$ws.Range("D23").Value = "'13.20"
$ws.Range("E23").Value = '  +1.86%  '

$ws.Range("D24").Value = "'27.396.62"
$ws.Range("E24").Value = '  +6.43%  '

$ws.Range("D25").Value = "'2.476"
$ws.Range("E25").Value = '  +0.46%  '

$ws.Range("D26").Value = "'2.736"
$ws.Range("E26").Value = '  +10.28%  '

$ws.Range("D27").Value = "'22.20"
$ws.Range("E27").Value = '  +16.06%  '

$ws.Range("E28").Value = '  +16.84%  '

$ws.Range("D29").Value = "'154.68"
$ws.Range("E29").Value = '  +3.50%  '

$ws.Range("D30").Value = "'2.021.35"
$ws.Range("E30").Value = '  +7.00%  '

$ws.Range("D31").Value = "'137.29"
$ws.Range("E31").Value = '  +7.23%  '

$ws.Range("D32").Value = "'6.429"
$ws.Range("E32").Value = '  +8.06%  '

$ws.Range("D33").Value = "'4.120"
$ws.Range("E33").Value = '  +0.09%  '

$ws.Range("D34").Value = "'13.90"
$ws.Range("E34").Value = '  +8.33%  '

$ws.Range("D35").Value = "'0.08847"
$ws.Range("E35").Value = '  +4.08%  '

$ws.Range("D36").Value = "'1.724"
$ws.Range("E36").Value = '  +2.97%  '

$ws.Range("D37").Value = "'5.680"
$ws.Range("E37").Value = '  +7.02%  '

$ws.Range("D38").Value = "'0.7113"
$ws.Range("E38").Value = '  +16.45%  '

$ws.Range("D39").Value = "'0.06591"
$ws.Range("E39").Value = '  +5.93%  '

$ws.Range("D40").Value = "'0.2275"
$ws.Range("E40").Value = '  +7.22%  '

$ws.Range("D41").Value = "'0.02423"
$ws.Range("E41").Value = '  +7.65%  '

$ws.Range("D42").Value = "'9.066"
$ws.Range("E42").Value = '  +6.42%  '

$ws.Range("D43").Value = "'1.261"
$ws.Range("E43").Value = '  -0.73%  '

$ws.Range("D44").Value = "'15.07"
$ws.Range("E44").Value = '  +4.21%  '

$ws.Range("D45").Value = "'0.6653"
$ws.Range("E45").Value = '  +14.12%  '

$ws.Range("E46").Value = '  +0.16%  '

$ws.Range("D47").Value = "'3.979"
$ws.Range("E47").Value = '  +3.45%  '

$ws.Range("D48").Value = "'2.200"
$ws.Range("E48").Value = '  +9.90%  '

$ws.Range("D49").Value = "'133.26"
$ws.Range("E49").Value = '  +5.61%  '

$ws.Range("D50").Value = "'0.07378"
$ws.Range("E50").Value = '  +2.20%  '

$ws.Range("D51").Value = "'81.11"
$ws.Range("E51").Value = '  +6.04%  '
